$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = "W"

$ws.Columns.Item(1).ColumnWidth = 22

$ws.Range("C8").Select()
